$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Insert a new row at position 10, pushing existing rows (10-25) down to (11-26)
$ws.Rows.Item(10).Insert()

# Fill in the new student record ("Rescatable") in the newly inserted row
$ws.Cells.Item(10, 1).Value = 18330051920213
$ws.Cells.Item(10, 2).Value = "LADINO"
$ws.Cells.Item(10, 3).Value = "URBINA"
$ws.Cells.Item(10, 4).Value = "MARIBEL"
$ws.Cells.Item(10, 5).Value = "INTRODUCCIÓN A LA ECONOMÍA"
$ws.Cells.Item(10, 6).Value = "6ARHM"
$ws.Cells.Item(10, 7).Value = 2
